$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows: row number, date serial (col A), col B, col C, col D
$newRows = @(
    @(465, 44539, 1, 12, 366.412213740458),
    @(466, 44540, 3, 15, 458.0152671755725),
    @(467, 44541, 0, 14, 427.4809160305344),
    @(468, 44542, 1, 12, 366.412213740458),
    @(469, 44543, 2, 10, 305.3435114503817),
    @(470, 44544, 0, 7, 213.7404580152672),
    @(471, 44545, 0, 7, 213.7404580152672),
    @(472, 44546, 2, 8, 244.2748091603054),
    @(473, 44547, 1, 6, 183.206106870229),
    @(474, 44548, 0, 6, 183.206106870229),
    @(475, 44550, 5, 10, 305.3435114503817),
    @(476, 44551, 1, 9, 274.8091603053435),
    @(477, 44552, 0, 9, 274.8091603053435),
    @(478, 44553, 2, 11, 335.8778625954199),
    @(479, 44554, 3, 12, 366.412213740458),
    @(480, 44555, 2, 13, 396.9465648854962),
    @(481, 44556, 3, 16, 488.5496183206107),
    @(482, 44557, 2, 13, 396.9465648854962),
    @(483, 44558, 0, 12, 366.412213740458),
    @(484, 44559, 2, 14, 427.4809160305344),
    @(485, 44560, 6, 18, 549.6183206106871),
    @(486, 44561, 2, 17, 519.0839694656488),
    @(487, 44562, 6, 21, 641.2213740458016),
    @(488, 44563, 5, 23, 702.2900763358779),
    @(489, 44564, 5, 26, 793.8931297709925),
    @(490, 44565, 0, 26, 793.8931297709925),
    @(491, 44566, 0, 24, 732.824427480916)
)

# Copy the formatting of the last existing data row (A464) down through
# the newly appended date cells (A465:A491) before writing values.
$ws.Range("A464").Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
